$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above row 5 for the new day's data (shifts rows 5:37 down to 6:38)
$ws.Rows("5:5").Insert()

# 2. Fix up the formatting of the newly inserted row 5 so it matches the rest of the
#    data rows (reuse existing styles instead of engine-generated ones) by copying
#    formats only from row 6 (the row that used to be row 5 before the insert).
$ws.Range("A6:G6").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the new row with the latest day's figures.
$ws.Range("A5").Value = 44336
$ws.Range("B5").Value = "(木)"
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value = 81422
$ws.Range("E5").Value = 109111

# 4. Update the running-total row (row 4) with the new cumulative figures.
$ws.Range("D4").Value = 3865493
$ws.Range("E4").Value = 2323873

# 5. Update the "as of" label shown above the table.
$ws.Range("E2").Value = "（5月20日時点）"
